$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column R (year 2021) picks up the formatting already used by column P
# (year 2019): border stays, vertical alignment switches from "top" to
# "center" (row 4/5), and the number format on row 6 drops the custom
# "0.0" format in favour of the plain format already used across L6:P6. ---
$ws.Range("R4").VerticalAlignment = -4108   # xlVAlignCenter
$ws.Range("R5").VerticalAlignment = -4108   # xlVAlignCenter

$ws.Range("P6").Copy()
$ws.Range("R6").PasteSpecial(-4122)         # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- New column S (year 2022), matching the (now updated) column R format ---

# Row 3: blank divider cell.
$ws.Range("R3").Copy()
$ws.Range("S3").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 4: header year.
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("S4").Value = 2022

# Row 5: first data series.
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("S5").Value = 1.8

# Row 6: second data series.
$ws.Range("R6").Copy()
$ws.Range("S6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("S6").Value = 8.4

# --- Selection moves to S3 ---
$ws.Range("S3").Select()
